$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.181.76'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '3.738.80'
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.02'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.91'
$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '662.52'
$ws.Range("E7").Value = '  +0.81%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.429'
$ws.Range("E8").Value = '  +1.65%  '

$ws.Range("E9").Value = '  -1.55%  '

$ws.Range("E10").Value = '  +0.00%  '

$ws.Range("D11").Value = '3.735.44'
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000324'
$ws.Range("E12").Value = '  +20.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.91'
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.209'
$ws.Range("E14").Value = '  +1.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.92'
$ws.Range("E15").Value = '  +1.08%  '

$ws.Range("D16").Value = '4.435.15'
$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("D17").Value = '97.089.56'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("E18").Value = '  +17.90%  '

$ws.Range("D19").Value = '3.730.90'
$ws.Range("E19").Value = '  +1.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.19'
$ws.Range("E20").Value = '  +3.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.87'
$ws.Range("E21").Value = '  -0.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.508'
$ws.Range("E22").Value = '  -3.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '529.74'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.50'
$ws.Range("E24").Value = '  +1.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000226'
$ws.Range("E25").Value = '  +10.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.93'
$ws.Range("E26").Value = '  -2.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '108.82'
$ws.Range("E27").Value = '  +6.77%  '

$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.194'
$ws.Range("E28").Value = '  +15.44%  '

$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.81'
$ws.Range("E29").Value = '  +4.49%  '

$ws.Range("D30").Value = '3.914.84'
$ws.Range("E30").Value = '  +1.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.04'
$ws.Range("E31").Value = '  +3.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.07'
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.192'
$ws.Range("E34").Value = '  +3.55%  '

$ws.Range("E35").Value = '  -2.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '33.21'
$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '652.05'
$ws.Range("E38").Value = '  -3.89%  '

$ws.Range("E39").Value = '  +0.85%  '

$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.169'
$ws.Range("E42").Value = '  +5.12%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.65'
$ws.Range("E43").Value = '  +4.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.89'
$ws.Range("E44").Value = '  +4.44%  '

$ws.Range("E45").Value = '  +2.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.990'
$ws.Range("E46").Value = '  +2.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.480'
$ws.Range("E47").Value = '  +8.62%  '

$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("E49").Value = '  +3.88%  '

$ws.Range("E50").Value = '  +1.91%  '
